$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A75 value (stored as text even though it looks numeric) ---
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "6585"
$ws.Range("A75").Style = "Normal"

# --- Add new row 87 data ---
# Text columns: force text entry (via temporary "@" format) so that
# numeric-looking strings / dates are preserved as plain text, matching
# the source data's inlineStr type, then restore the default "Normal"
# style so no extra formatting is introduced.
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "-549"
$ws.Range("A87").Style = "Normal"

$ws.Range("B87").NumberFormat = "@"
$ws.Range("B87").Value = "8/7/2025"
$ws.Range("B87").Style = "Normal"

$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Value = "14 de Julio 65"
$ws.Range("C87").Style = "Normal"

$ws.Range("D87").NumberFormat = "@"
$ws.Range("D87").Value = "13"
$ws.Range("D87").Style = "Normal"

$ws.Range("E87").NumberFormat = "@"
$ws.Range("E87").Value = "808749189"
$ws.Range("E87").Style = "Normal"

$ws.Range("F87").NumberFormat = "@"
$ws.Range("F87").Value = "AYKO"
$ws.Range("F87").Style = "Normal"

$ws.Range("G87").NumberFormat = "@"
$ws.Range("G87").Value = "Pendiente"
$ws.Range("G87").Style = "Normal"

$ws.Range("H87").NumberFormat = "@"
$ws.Range("H87").Value = "picada"
$ws.Range("H87").Style = "Normal"

# Numeric columns
$ws.Range("I87").Value = 1

$ws.Range("J87").NumberFormat = "@"
$ws.Range("J87").Value = "Cambio"
$ws.Range("J87").Style = "Normal"

$ws.Range("K87").NumberFormat = "@"
$ws.Range("K87").Value = "Sin equipos"
$ws.Range("K87").Style = "Normal"

$ws.Range("L87").NumberFormat = "@"
$ws.Range("L87").Value = "Pasante"
$ws.Range("L87").Style = "Normal"

$ws.Range("M87").Value = -58.468496
$ws.Range("N87").Value = -34.591282

$ws.Range("O87").NumberFormat = "@"
$ws.Range("O87").Value = "Paternal"
$ws.Range("O87").Style = "Normal"

$ws.Range("P87").NumberFormat = "@"
$ws.Range("P87").Value = "Capital Norte"
$ws.Range("P87").Style = "Normal"
